$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for the Price/Volume columns so numeric-looking
# strings (e.g. "0.9993", "244.00") are written as literal text, matching
# the source data (which also contains non-numeric strings like "30.082.38").
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.082.38"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.877.96"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "244.00"
$ws.Range("E5").Value = "  -2.46%  "
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").Value = "0.4913"
$ws.Range("E7").Value = "  -1.45%  "
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("D9").Value = "0.06621"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "1.879.89"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "16.64"
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("D12").Value = "0.07203"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "0.6686"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("D14").Value = "86.47"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "4.918"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "30.035.36"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("D18").Value = "0.9993"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "12.81"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "2.122.69"
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("D21").Value = "0.9983"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "4.792"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").Value = "5.883"
$ws.Range("E23").Value = "  +3.66%  "
$ws.Range("D24").Value = "9.142"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("D25").Value = "152.61"
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("D26").Value = "143.31"
$ws.Range("E26").Value = "  +8.29%  "
$ws.Range("D27").Value = "16.98"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("D28").Value = "1.900"
$ws.Range("E28").Value = "  -3.30%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "4.203"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "0.08771"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "3.997"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").Value = "0.05062"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").Value = "0.7220"
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").Value = "1.114"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").Value = "2.660"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("D37").Value = "0.01855"
$ws.Range("E37").Value = "  +11.53%  "
$ws.Range("D38").Value = "2.691"
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("D39").Value = "2.160"
$ws.Range("E39").Value = "  -3.92%  "
$ws.Range("D40").Value = "0.9319"
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("D41").Value = "5.780"
$ws.Range("E41").Value = "  -3.86%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.4237"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "0.9990"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").Value = "103.29"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "7.393"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").Value = "0.1276"
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("D47").Value = "0.05713"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").Value = "32.86"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "0.3788"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "8.272"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "1.346"
$ws.Range("E51").Value = "  +0.17%  "

# Restore original (unformatted) style on the touched range so the only
# change recorded is the literal value, not the number format.
$priceRange.ClearFormats()
